# Block all offshore wind not currently planned/under construction in Current Policies
#
# 1) BBNPPTY sheet: offshore wind (row 15) years 2024-2050 (cols E:AE) flip
#    from 0 (banned / blocked, integer-formatted style) to 1 (allowed,
#    default/general format) -- i.e. all years now allowed, matching the
#    "already planned capacity" policy captured elsewhere.
# 2) About sheet: append three new note lines (rows 18-20) explaining the
#    offshore wind assumption, which also adds matching shared strings.
# 3) Leave "About" as the sheet the user ends up looking at (tab activated,
#    selection sitting just below the new notes), and leave BBNPPTY's
#    selection/scroll parked on the row that was just edited.

$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("BBNPPTY")

# --- 1) Un-ban offshore wind (row 15) for every forecast year, 2024-2050 ---
$offshoreRange = $wsData.Range("E15:AE15")
$offshoreRange.Value = 1
# The "banned" years previously carried an integer number-format style;
# clear it back to Normal/General now that they're simple 0/1 flags like
# the rest of the row.
$offshoreRange.Style = "Normal"

# --- 2) Add explanatory notes under the existing About-sheet notes ---
$wsAbout.Range("A18").Value = "Due to recent blocks on offshore wind permitting and leasing, we assume only"
$wsAbout.Range("A19").Value = "already planned capacity (captured in elec/BPMCCS) will be built, and other"
$wsAbout.Range("A20").Value = "economic or reliability additions will not occur (offshore wind set to 1 in this file)."

# --- 3) View-state: leave the About sheet active/selected, with the data ---
#        sheet scrolled/selected to the edited row for reference.
$wsData.Activate()
$wsData.Range("D15:AE15").Select()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1

$wsAbout.Activate()
$wsAbout.Range("A21").Select()
